$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for data rows 2-18
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05)
for ($row = 2; $row -le 18; $row++) {
    $ws.Range("C$row").Value = 45204
}
